# Replace the reviewer's name in the signature line, and move the hidden
# "_GoBack" bookmark from the QR-code picture (where it happened to sit)
# to the spot that was just edited - exactly what Word itself does after
# you type over a selection.

$d = $word.ActiveDocument

# 1. Drop the bookmark from its old location around the QR-code image.
$d.Bookmarks("_GoBack").Delete()

# 2. Locate the run that holds the old name.
$target = $d.Content
$target.Find.Execute("Иванова Е.А.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$rStart = $target.Start
$rEnd = $target.End

# 3. Pin both edges of that run with throw-away bookmarks *before* touching
#    the text - this stops the engine from silently fusing the edited run
#    back into its neighbours ("(" .. ")") when the paragraph is rebuilt.
$d.Bookmarks.Add("zzTmpBefore", $d.Range($rStart, $rStart)) | Out-Null
$d.Bookmarks.Add("zzTmpAfter", $d.Range($rEnd, $rEnd)) | Out-Null

# 4. Swap in the new name.
$d.Range($rStart, $rEnd).Text = "Киселев В. С."

# 5. Discard the leading pin, and turn the trailing pin into the real
#    "_GoBack" bookmark, which is where Word leaves it after an edit.
$d.Bookmarks("zzTmpBefore").Delete()
$trailing = $d.Bookmarks("zzTmpAfter")
$goBackRange = $trailing.Range
$trailing.Delete()
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
